$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Alegria)
$ws.Range("D2").Value = 7.95
$ws.Range("E2").Value = 16.760000000000002
$ws.Range("G2").Value = 8.2330000000000005
$ws.Range("H2").Value = 7.2690000000000001
$ws.Range("I2").Value = 20.53
$ws.Range("J2").Value = 11.69

# Row 3 (Bodega Bay)
$ws.Range("D3").Value = 7.8940000000000001
$ws.Range("E3").Value = 13.46
$ws.Range("G3").Value = 8.532
$ws.Range("H3").Value = 7.5010000000000003
$ws.Range("I3").Value = 22.239000000000001
$ws.Range("J3").Value = 6.3330000000000002

# Row 4 (Lompoc Landing)
$ws.Range("D4").Value = 7.8490000000000002
$ws.Range("E4").Value = 14.63
$ws.Range("G4").Value = 8.2490000000000006
$ws.Range("H4").Value = 7.1790000000000003
$ws.Range("I4").Value = 22.88
$ws.Range("J4").Value = 11.07

# Update the active cell selection to match the saved view state
$ws.Range("H12").Select()
